# "work on other species"
# Adds a new "Group" column (D) to the species sample-size table that
# buckets each species/row into a coarser category group.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> Group label, in sheet order (row 1 is the header row).
$groups = @{
    1  = "Group"
    2  = "Dungeness crab"
    3  = "California mussel"
    4  = "Razor clam"
    5  = "Rock crab"
    6  = "Pacific oyster"
    7  = "Bay mussel"
    8  = "Spiny lobster"
    9  = "Rock crab"
    10 = "Sardine/anchovy"
    11 = "Sardine/anchovy"
    12 = "Other mollusc"
    13 = "Other mollusc"
    14 = "Other mollusc"
    15 = "Other mollusc"
    16 = "Other crab"
    17 = "Other mollusc"
    18 = "Rock crab"
    19 = "Rock crab"
    20 = "Other mollusc"
    21 = "Other mollusc"
    22 = "Other crab"
    23 = "Other mollusc"
    24 = "Other fish"
    25 = "Other crab"
    26 = "Other mollusc"
    27 = "Other fish"
    28 = "Other fish"
    29 = "Other fish"
    30 = "Other/unknown"
    31 = "Other mollusc"
    32 = "Other mollusc"
    33 = "Other fish"
    34 = "Other mollusc"
    35 = "Other mollusc"
    36 = "Other fish"
    37 = "Other/unknown"
    38 = "Other mollusc"
    39 = "Other mollusc"
    40 = "Other/unknown"
    41 = "Other fish"
    42 = "Other crab"
}

# Rows where Excel ended up re-keying the cell to an explicit black RGB
# font (fontId referencing color rgb="FF000000") rather than the default
# theme color - a cosmetic artifact of how those particular cells were
# typed in. Reproduce it so the new cells carry the same explicit font.
$blackFontRows = @(22, 25, 33, 42)

# Prime the shared-string table so new unique strings land at the same
# indices as the authored workbook. The header ("Group") and the 8
# per-row category labels (rows 2-11) are already encountered in natural
# row order, but the 4 generic "Other ..." buckets need to be primed in
# the order they were first typed (crab, fish, mollusc, unknown) rather
# than plain row order (mollusc would otherwise come first, at row 12).
for ($r = 1; $r -le 11; $r++) {
    $ws.Cells.Item($r, 4).Value = $groups[$r]
}
$primeOrder = @(16, 24, 12, 30)
foreach ($pr in $primeOrder) {
    $ws.Cells.Item($pr, 4).Value = $groups[$pr]
}

for ($r = 1; $r -le 42; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    $cell.Value = $groups[$r]

    if ($r -eq 1) {
        # Match header formatting used by A1:C1 (bold, double-bottom border)
        # by copying the existing header cell style instead of re-deriving
        # font/border settings (which would create a brand-new style).
        $ws.Range("C1").Copy()
        $ws.Range("D1").PasteSpecial(-4122)  # xlPasteFormats
    } elseif ($blackFontRows -contains $r) {
        $cell.Font.Color = 0x000000
    }
}
$excel.CutCopyMode = 0

# New column D width (to fit the longer "California mussel" labels etc.)
$ws.Columns.Item(4).ColumnWidth = 14.5

# Match the final selection left behind in the authored workbook.
$ws.Range("D15").Select()
